$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row: add C2 "measuring, min" first, then B2 "Δt, ms",
# so the shared-string table indices come out in the same order as the target file
$ws.Range("C2").Value = "measuring, min"
$ws.Range("B2").Value = "Δt, ms"

# Set column C width (~16.43 chars, matches autofit width for "measuring, min")
$ws.Columns.Item(3).ColumnWidth = 15.7114955357143

# Fill in data rows 3-6 (distance, Δt(ms), measuring(min))
$ws.Range("A3").Value = 0.5
$ws.Range("B3").Value = 18
$ws.Range("C3").Value = 1

$ws.Range("A4").Value = 1
$ws.Range("B4").Value = 18
$ws.Range("C4").Value = 1

$ws.Range("A5").Value = 2
$ws.Range("B5").Value = 18
$ws.Range("C5").Value = 1

$ws.Range("A6").Value = 7
$ws.Range("B6").Value = 18
$ws.Range("C6").Value = 1

# Remove old row 7 (distance value 10) entirely
$ws.Range("A7").ClearContents()

# Add conclusion row
$ws.Range("A8").Value = "Conclusion: Δt does not depend on the distance between the modules"

# Update selection to match target workbook state
$ws.Range("C16").Select()
